# Apply the "indian foods, rm sushi" edit to runs.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 (fast_food) restaurant names updated
$ws.Range("A11").Value = "Dave's Sandwiches,Burger Masters,Metro Drive-In"

# Row 9 was the sushi row -> replace with the new indian_food row
$ws.Range("B9").Value = "indian_food"
$ws.Range("A9").Value = "Little Delhi House,Curry Leaf Restaurant,Tandoori Kitchen"

# Widen column C a bit (matches new col width in the diff)
$ws.Columns.Item(3).ColumnWidth = 25.33203125

# Move the active selection to A10 (as in the edited file)
$ws.Range("A10").Select()
